$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the PDF paths in column K for rows 2-13: the folder "PDF" was
#    duplicated in the path (...\PDF\file.pdf -> ...\PDF\PDF\file.pdf).
$ws.Range("K2").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia agrol  diciembre 2022.pdf"
$ws.Range("K3").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia bricks febrero 2023.pdf"
$ws.Range("K4").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia codigo a enero 2023.pdf"
$ws.Range("K5").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia cuadras asesores diciembre 2022.pdf"
$ws.Range("K6").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia desoflex 19-03-2020.pdf"
$ws.Range("K7").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia diseños luna 10-01-23.pdf"
$ws.Range("K8").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\Constancia Interexporta MARZO 2023.pdf"
$ws.Range("K9").Value  = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia jose luis zarate losa 19-01-23.pdf"
$ws.Range("K10").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia moises mercado torres 03-11-22.pdf"
$ws.Range("K11").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia monica alejandra zarate losa  22-02-23.pdf"
$ws.Range("K12").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia pascual ibarra 23-01-23.pdf"
$ws.Range("K13").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia rigoberto mora 10-01-23.pdf"

# 2) Insert a new row at position 14 - this pushes the former row 14
#    (JAUREGUI / rosa karmin jauregui) down to row 15, and updates the
#    used-range dimension to A1:K15 automatically.
$ws.Rows("14").Insert()

# 3) The row that was pushed down to row 15 also needs its PDF path fixed
#    the same way as the other rows.
$ws.Range("K15").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia rosa karmin jauregui 16-02-23.pdf"

# 4) Populate the new row 14 with the "MOISES/RIGOBERTO MORA ALVAREZ"
#    (reactivado) record and its newly generated PDF.
$ws.Range("A14").Value = "Constancia"
$ws.Range("B14").Value = "MOAR741018D36 "
$ws.Range("C14").Value = " MOAR741018HJCRLG02 "
$ws.Range("D14").Value = " MOAR741018HJCRLG02 "
$ws.Range("E14").Value = " MORA "
$ws.Range("F14").Value = " ALVAREZ "
$ws.Range("H14").Value = " 01DEJUNIODE1996"
$ws.Range("I14").Value = " REACTIVADO "
$ws.Range("J14").Value = "Regímenes:   Régimen Fecha Inicio Fecha Fin Régimen de Ingresos por Dividendos (socios y accionistas) 01/01/2017 Régimen de las Personas Físicas con Actividades Empresariales y Profesionales 01/01/2019 Obligaciones"
$ws.Range("K14").Value = "C:\Users\victo\Documents\Proyectos\Automatizacion-Cuadras\PDF\PDF\constancia rigoberto mora ENERO 2023.pdf"
